$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 66

$ws.Cells.Item($newRow, 1).Value = "2025-08-28 03:47:45 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-28 09:17:45 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

# Match the style of the row above (row 65) for every cell in the new row
$ws.Range("A65:H65").Copy()
$ws.Range("A66:H66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
